$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "248.53", "1.000"); Excel
# auto-converts such strings to numbers on assignment. Toggle the cell to Text
# format for the assignment, then back to General, so the stored type/value
# both match the source data (plain text, trailing zeros preserved).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.759.43"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.76"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.53"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4737"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2932"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06541"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.08"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07806"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.18"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.891.29"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7364"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.256"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "283.13"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.891.96"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.21"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007546"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.139.43"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.335"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.268"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.249"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.34"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.95"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.927"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.342"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09728"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.500"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.306"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.204"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04864"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.128"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6990"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01902"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.808"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.380"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "76.23"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.021"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4270"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8354"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.31"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.497"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.71"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.055"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "920.09"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05760"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +2.15%  "
